# Rename the wage & contribution column headers on both the "Jack" and
# "Jill" sheets for consistency, and bold the header row (selecting the
# whole row first, matching how this would be done interactively).
#
# Column label changes (same on both sheets):
#   C1  "ctrb taxable"   -> "taxable ctrb"
#   D1  "ctrb 401k"      -> "401k ctrb"
#   E1  "ctrb Roth 401k" -> "Roth 401k ctrb"
#   F1  "ctrb IRA"       -> "IRA ctrb"
#   G1  "ctrb Roth IRA"  -> "Roth IRA ctrb"
#   H1  "Roth X"         -> "Roth conv"

$wb = $excel.ActiveWorkbook

$sheetNames = @("Jill", "Jack")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("H1").Value = "Roth conv"
    $ws.Range("C1").Value = "taxable ctrb"
    $ws.Range("D1").Value = "401k ctrb"
    $ws.Range("E1").Value = "Roth 401k ctrb"
    $ws.Range("F1").Value = "IRA ctrb"
    $ws.Range("G1").Value = "Roth IRA ctrb"

    # Select the header row and bold it, leaving the row selected
    # (Jack is processed last so it ends up the active sheet/selection).
    $ws.Rows("1:1").Select() | Out-Null
    $ws.Rows("1:1").Font.Bold = $true
}
